$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D - total days
    $eCell = $ws.Cells.Item($r, 5)   # column E - remaining days
    $fCell = $ws.Cells.Item($r, 6)   # column F - start date (yyyymmdd)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null) { continue }

    # Skip rows whose start date is not a well-formed yyyymmdd date
    # (mirrors the source automation's behaviour, which errors out and
    # leaves such rows untouched - e.g. a corrupted date like 202510929).
    $fText = [string][int]$fVal
    $isValidDate = $false
    if ($fText.Length -eq 8) {
        $y = [int]$fText.Substring(0, 4)
        $mo = [int]$fText.Substring(4, 2)
        $da = [int]$fText.Substring(6, 2)
        if ($mo -ge 1 -and $mo -le 12 -and $da -ge 1 -and $da -le 31) {
            try {
                $null = Get-Date -Year $y -Month $mo -Day $da
                $isValidDate = $true
            } catch {
                $isValidDate = $false
            }
        }
    }

    if (-not $isValidDate) { continue }

    if ($eVal -le 1) {
        # Cycle finished: reset to a fresh full cycle starting "today" (2026-02-07)
        $dVal = $dCell.Value2
        $eCell.Value2 = $dVal
        $fCell.Value2 = 20260207
    } else {
        # One more day has elapsed since the last update
        $eCell.Value2 = $eVal - 1
    }
}
